$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8587174348697395
$ws.Range("B3").Value = 0.8217821782178217
$ws.Range("B4").Value = 0.8680904522613065
$ws.Range("B5").Value = 0.6125461254612546
$ws.Range("B6").Value = 0.9504814305364512
$ws.Range("B7").Value = 0.9229470371660282
